{"js": "// Add a short \"body\" paragraph (style \"First Paragraph\") right after each\n// section heading in the issue template, per the commit's docx conversion\n// of the markdown issue template.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Heading text -> new paragraph text to insert right after it.\n// (The \"Environment\" heading is intentionally skipped - the diff does not\n// add a paragraph there, it already has its own bullet list.)\nconst insertions = [\n  { after: \"Issue Template\", text: \"Release: v1.0.0\", bold: true },\n  { after: \"Description\", text: \"Describe the issue or feature request in detail.\" },\n  { after: \"Steps to Reproduce\", text: \"List steps to reproduce the issue, if applicable.\" },\n  { after: \"Expected Behavior\", text: \"What did you expect to happen?\" },\n  { after: \"Actual Behavior\", text: \"What actually happened?\" },\n  { after: \"Additional Context\", text: \"Add any other context or screenshots about the issue here.\" },\n];\n\nfunction findParagraphByText(targetText) {\n  for (const p of paragraphs.items) {\n    if (p.text === targetText) return p;\n  }\n  throw new Error(\"Could not locate paragraph with text: \" + targetText);\n}\n\nfor (const insertion of insertions) {\n  const anchor = findParagraphByText(insertion.after);\n  const newParagraph = anchor.insertParagraph(insertion.text, Word.InsertLocation.after);\n  newParagraph.style = \"First Paragraph\";\n  if (insertion.bold) {\n    newParagraph.font.bold = true;\n    // Mirror Word's usual behavior of also flagging bold for complex-script\n    // runs (<w:bCs/>) alongside the regular bold (<w:b/>).\n    newParagraph.font.boldBidirectional = true;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Add a short \"body\" paragraph (style \"First Paragraph\") right after each\n# section heading in the issue template, per the commit's docx conversion\n# of the markdown issue template.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphIndexByText($doc, $text) {\n    $i = 1\n    foreach ($p in $doc.Paragraphs) {\n        $t = $p.Range.Text.TrimEnd(\"`r\", \"`n\")\n        if ($t -eq $text) {\n            return $i\n        }\n        $i++\n    }\n    return -1\n}\n\nfunction Insert-BodyParagraphAfter($doc, $headingText, $newText, $bold) {\n    $idx = Get-ParagraphIndexByText $doc $headingText\n    $heading = $doc.Paragraphs($idx)\n\n    # Insert a new (blank) paragraph right after the heading, then grab it\n    # by its now-known index so we don't depend on Range.Next() semantics.\n    $heading.Range.InsertParagraphAfter()\n    $newPara = $doc.Paragraphs($idx + 1)\n    $newPara.Style = \"First Paragraph\"\n    $newPara.Range.Text = $newText\n\n    # Only format the text itself, not the trailing paragraph mark, so the\n    # bold doesn't leak into the paragraph-mark run properties.\n    $r = $newPara.Range\n    $r.MoveEnd(1, -1)\n    if ($bold) {\n        $r.Font.Bold = 1\n        $r.Font.BoldBi = 1\n    }\n}\n\n# \"Environment\" is intentionally skipped - the diff does not add a\n# paragraph there, it already has its own bullet list.\nInsert-BodyParagraphAfter $d \"Issue Template\" \"Release: v1.0.0\" $true\nInsert-BodyParagraphAfter $d \"Description\" \"Describe the issue or feature request in detail.\" $false\nInsert-BodyParagraphAfter $d \"Steps to Reproduce\" \"List steps to reproduce the issue, if applicable.\" $false\nInsert-BodyParagraphAfter $d \"Expected Behavior\" \"What did you expect to happen?\" $false\nInsert-BodyParagraphAfter $d \"Actual Behavior\" \"What actually happened?\" $false\nInsert-BodyParagraphAfter $d \"Additional Context\" \"Add any other context or screenshots about the issue here.\" $false\n"}
